# Update countries & provincias Spain
# Refreshes the COVID-19 case-count snapshot on the "Pais" sheet:
#   - bumps the "Datos actualizados ..." timestamp in A1 from 18:22 to 18:52
#   - updates the per-country Casos totales / Nuevos casos / Casos activos /
#     Recuperados / Casos criticos / Muertes hoy / Muertes figures (columns B-H)
#     for every country whose numbers moved between the two snapshots.
# Country names in column A keep their existing row positions - only the
# numeric statistics and the timestamp text change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp cell
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 18:52"

# Estados Unidos
$ws.Range("B4").Value = 293494
$ws.Range("C4").Value = 16333
$ws.Range("D4").Value = 14436
$ws.Range("E4").Value = 271162
$ws.Range("F4").Value = 6219
$ws.Range("G4").Value = 492
$ws.Range("H4").Value = 7896

# Chequia
$ws.Range("B26").Value = 4604
$ws.Range("C26").Value = 331
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = 4442
$ws.Range("F26").Value = 148
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 137

# Irlanda
$ws.Range("B27").Value = 4362
$ws.Range("C27").Value = 172
$ws.Range("D27").Value = 78
$ws.Range("E27").Value = 4225
$ws.Range("F27").Value = 87
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 59

# Ecuador
$ws.Range("B33").Value = 3465
$ws.Range("C33").Value = 97
$ws.Range("E33").Value = 3228
$ws.Range("G33").Value = 27
$ws.Range("H33").Value = 172

# Luxemburgo
$ws.Range("B38").Value = 2729
$ws.Range("C38").Value = 117
$ws.Range("E38").Value = 2198
$ws.Range("F38").Value = 35

# Colombia
$ws.Range("B53").Value = 1325
$ws.Range("C53").Value = 250
$ws.Range("D53").Value = 109
$ws.Range("E53").Value = 1213
$ws.Range("F53").Value = 37
$ws.Range("H53").Value = 3

# Argelia
$ws.Range("B54").Value = 1267
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 55
$ws.Range("E54").Value = 1187
$ws.Range("F54").Value = 50
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 25

# Singapur
$ws.Range("B55").Value = 1251
$ws.Range("C55").Value = 80
$ws.Range("D55").Value = 90
$ws.Range("E55").Value = 1031
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 25
$ws.Range("H55").Value = 130

# Croacia
$ws.Range("B56").Value = 1189
$ws.Range("C56").Value = 75
$ws.Range("D56").Value = 297
$ws.Range("E56").Value = 886
$ws.Range("F56").Value = 24
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 6

# Ucrania
$ws.Range("B57").Value = 1126
$ws.Range("C57").Value = 47
$ws.Range("D57").Value = 119
$ws.Range("E57").Value = 995
$ws.Range("F57").Value = 39
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 12

# Catar
$ws.Range("B58").Value = 1096
$ws.Range("C58").Value = 24
$ws.Range("D58").Value = 23
$ws.Range("E58").Value = 1045
$ws.Range("F58").Value = 16
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 28

# Barein
$ws.Range("D70").Value = 423
$ws.Range("E70").Value = 261

# Albania
$ws.Range("B89").Value = 334
$ws.Range("C89").Value = 13
$ws.Range("D89").Value = 40
$ws.Range("E89").Value = 294
$ws.Range("F89").Value = 4
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 0

# Reunion
$ws.Range("B90").Value = 333
$ws.Range("C90").Value = 29
$ws.Range("D90").Value = 99
$ws.Range("E90").Value = 216
$ws.Range("F90").Value = 7
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 18

# Honduras
$ws.Range("B96").Value = 266
$ws.Range("C96").Value = 39
$ws.Range("D96").Value = 25
$ws.Range("E96").Value = 239
$ws.Range("F96").Value = 8
$ws.Range("H96").Value = 2

# Uzbekistan
$ws.Range("B97").Value = 264
$ws.Range("C97").Value = 42
$ws.Range("D97").Value = 3
$ws.Range("E97").Value = 246
$ws.Range("F97").Value = 10
$ws.Range("H97").Value = 15

# Montenegro
$ws.Range("B106").Value = 201
$ws.Range("C106").Value = 27
$ws.Range("E106").Value = 198

# Georgia
$ws.Range("B110").Value = 159
$ws.Range("C110").Value = 4
$ws.Range("E110").Value = 127

# Guadalupe
$ws.Range("B117").Value = 134
$ws.Range("C117").Value = 6
$ws.Range("D117").Value = 14
$ws.Range("E117").Value = 118
$ws.Range("F117").Value = 3
$ws.Range("H117").Value = 2

# Mayotte
$ws.Range("B118").Value = 130
$ws.Range("D118").Value = 24
$ws.Range("E118").Value = 99
$ws.Range("F118").Value = 14
$ws.Range("H118").Value = 7

# Monaco
$ws.Range("C131").Value = 2
$ws.Range("D131").Value = 1
$ws.Range("E131").Value = 63
$ws.Range("F131").Value = 0
$ws.Range("H131").Value = 0

# Aruba
$ws.Range("B132").Value = 64
$ws.Range("D132").Value = 3
$ws.Range("E132").Value = 60
$ws.Range("F132").Value = 2
$ws.Range("H132").Value = 1
